$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -10520.89681921756
$ws.Range("C2").Value = 21184.44377802947
$ws.Range("D2").Value = -10000.49245589557
$ws.Range("E2").Value = -663.0545029160467
$ws.Range("F2").Value = 61.93692125772179
$ws.Range("G2").Value = 45.15013065477399
$ws.Range("H2").Value = 62.25406853330918
$ws.Range("I2").Value = 45.3534740409537
$ws.Range("J2").Value = 61.93692125784037
$ws.Range("K2").Value = 45.98047043381223
$ws.Range("L2").Value = 55.01990173258815
$ws.Range("O2").Value = 49.86825444426466
$ws.Range("P2").Value = 59.28641779413653
$ws.Range("R2").Value = 8.953380639800754
$ws.Range("S2").Value = -17.90676128295929
$ws.Range("T2").Value = 8.953380643158548
$ws.Range("X2").Value = -92.33646748703023
$ws.Range("Y2").Value = -142.3360080380841
$ws.Range("Z2").Value = -92.33646748703015
$ws.Range("AE2").Value = -8.953380639800754
$ws.Range("AF2").Value = 8.953380643158548
$ws.Range("AG2").Value = 8.953380639800754
$ws.Range("AH2").Value = -17.90676128295929
$ws.Range("AI2").Value = 8.953380643158548
$ws.Range("AJ2").Value = 8.953380639800754
$ws.Range("AK2").Value = -8.953380643158548
$ws.Range("AL2").Value = 24.99977027552693
$ws.Range("AM2").Value = -24.99977027552692
$ws.Range("AN2").Value = -92.33646748703023
$ws.Range("AO2").Value = -142.3360080380841
$ws.Range("AP2").Value = -92.33646748703015
$ws.Range("AQ2").Value = -24.99977027552693
$ws.Range("AR2").Value = 24.99977027552692
$ws.Range("AS2").Value = 61.93692125772179
$ws.Range("AT2").Value = 61.93692125772179
$ws.Range("AU2").Value = 62.25406853330901
$ws.Range("AV2").Value = 62.25406853330901
$ws.Range("AW2").Value = 62.25406853330918
$ws.Range("AX2").Value = 61.93692125784037
$ws.Range("AY2").Value = 61.93692125784037
$ws.Range("AZ2").Value = 45.15013065477399
$ws.Range("BA2").Value = 45.15013065477399
$ws.Range("BB2").Value = 45.3534740409537
$ws.Range("BC2").Value = 44.94092423767461
$ws.Range("BD2").Value = 45.76602384407818
$ws.Range("BE2").Value = 45.98047043381223
$ws.Range("BF2").Value = 45.98047043381217
$ws.Range("BG2").Value = 55.01990173258815
$ws.Range("BJ2").Value = 49.86825444426466
$ws.Range("BK2").Value = 59.28641779413653
